# Adds the "Introduction ro VBA" video-link row to the learning-videos sheet.
# A new row is inserted above the current row 17 (the first row of the
# Module-4-PQ block), pushing that block (and everything after it) down by
# one row, and the new row is filled in with the same layout used by the
# other "Module-3-XL" sub-topic rows immediately above it (rows 15 & 16):
# Module name in column B, topic text in column C, and a hyperlinked video
# URL in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 17 - everything currently on/under row 17
# (rows 17-64) shifts down to rows 18-65.
$ws.Rows.Item(17).Insert()

# Match the formatting used by the sibling rows above (15 & 16): body font
# for the Module/Topic columns, and the existing Hyperlink-column style for
# the link column.
$ws.Range("B17:C17").Font.Name = "Avenir Next LT Pro"
$ws.Range("B17:C17").Font.Size = 11

$ws.Range("B17").Value = "Module-3-XL"
$ws.Range("C17").Value = "Introduction ro VBA"

$videoUrl = "https://mentorskool-platform-uploads.s3.ap-south-1.amazonaws.com/strapiUploads/imageAssets/vba_basic_script_475f65923e.mp4"
$ws.Range("D17").Value = $videoUrl
$ws.Hyperlinks.Add($ws.Range("D17"), $videoUrl)
$ws.Range("D17").Style = $ws.Range("D16").Style

$ws.Rows.Item(17).RowHeight = 15

# Restore selection near the edited area, matching where the author was
# working after making the change.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A15").Select()
